# Auto-generated edit script applying numeric corrections to H/I/J/K/L/M/N
# columns across several rows in the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 489.53845
$ws.Range("I18").Value = 489.53845
$ws.Range("K18").Value = 489.53845
$ws.Range("M18").Value = -205.53845
# Row 43
$ws.Range("H43").Value = 3729.8
$ws.Range("I43").Value = 3166.3333
$ws.Range("K43").Value = 3166.3333
$ws.Range("M43").Value = -3097.3333
# Row 98
$ws.Range("H98").Value = 57981.934
$ws.Range("I98").Value = 62081.07
$ws.Range("K98").Value = 62081.07
$ws.Range("M98").Value = -60583.07
# Row 122
$ws.Range("H122").Value = 57981.934
$ws.Range("I122").Value = 62081.07
$ws.Range("K122").Value = 186243.21
$ws.Range("M122").Value = -183793.21
# Row 129
$ws.Range("H129").Value = 15461.134
$ws.Range("I129").Value = 771.5
$ws.Range("K129").Value = 2314.5
$ws.Range("M129").Value = 2685.5
# Row 137
$ws.Range("H137").Value = 32788758
$ws.Range("I137").Value = 18869822
$ws.Range("K137").Value = 56609466
$ws.Range("M137").Value = -56606916
# Row 138
$ws.Range("H138").Value = 4019007.8
$ws.Range("J138").Value = 5053791.5
$ws.Range("L138").Value = 15161374.5
$ws.Range("N138").Value = -15171654.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 16407043
$ws.Range("I32").Value = 27783992
$ws.Range("K32").Value = 27783992
$ws.Range("M32").Value = -27783705
# Row 45
$ws.Range("H45").Value = 1704.7142
$ws.Range("I45").Value = 1496.6666
$ws.Range("K45").Value = 1496.6666
$ws.Range("M45").Value = -1119.6666
# Row 61
$ws.Range("H61").Value = 27780232
$ws.Range("I61").Value = 35715724
$ws.Range("K61").Value = 35715724
$ws.Range("M61").Value = -35715512
# Row 74
$ws.Range("H74").Value = 47674796
$ws.Range("I74").Value = 50058036
$ws.Range("K74").Value = 50058036
$ws.Range("M74").Value = -50057162
# Row 77
$ws.Range("H77").Value = 47674796
$ws.Range("I77").Value = 50058036
$ws.Range("K77").Value = 250290180
$ws.Range("M77").Value = -250285812
# Row 110
$ws.Range("H110").Value = 8988.816000000001
$ws.Range("I110").Value = 9648.929
$ws.Range("K110").Value = 9648.929
$ws.Range("M110").Value = -7603.929
# Row 132
$ws.Range("H132").Value = 18523912
$ws.Range("I132").Value = 5785.3267
$ws.Range("K132").Value = 17355.9801
$ws.Range("M132").Value = -14825.9801
# Row 136
$ws.Range("H136").Value = 27780232
$ws.Range("I136").Value = 35715724
$ws.Range("K136").Value = 107147172
$ws.Range("M136").Value = -107144622

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 8
$ws.Range("H8").Value = 20147
$ws.Range("J8").Value = 7444
$ws.Range("L8").Value = 7444
$ws.Range("N8").Value = -7724
# Row 94
$ws.Range("H94").Value = 2299.625
$ws.Range("I94").Value = 2066.1667
$ws.Range("K94").Value = 2066.1667
$ws.Range("M94").Value = -1615.1667
# Row 105
$ws.Range("H105").Value = 13676
$ws.Range("I105").Value = 21031.6
$ws.Range("K105").Value = 21031.6
$ws.Range("M105").Value = -19284.6
# Row 107
$ws.Range("H107").Value = 6079.8
$ws.Range("I107").Value = 5634.3335
$ws.Range("K107").Value = 5634.3335
$ws.Range("M107").Value = -3714.3335
# Row 110
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
# Row 134
$ws.Range("H134").Value = 3828.6345
$ws.Range("I134").Value = 3380.4883
$ws.Range("K134").Value = 10141.4649
$ws.Range("M134").Value = -7606.464899999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 86
$ws.Range("H86").Value = 3753.8462
$ws.Range("I86").Value = 2843.7144
$ws.Range("K86").Value = 2843.7144
$ws.Range("M86").Value = -1720.7144
# Row 89
$ws.Range("H89").Value = 3753.8462
$ws.Range("I89").Value = 2843.7144
$ws.Range("K89").Value = 14218.572
$ws.Range("M89").Value = -8602.572
# Row 107
$ws.Range("H107").Value = 2110.4375
$ws.Range("I107").Value = 1446.6666
$ws.Range("J107").Value = 2963.8572
$ws.Range("K107").Value = 1446.6666
$ws.Range("L107").Value = 2963.8572
$ws.Range("M107").Value = 473.3334
$ws.Range("N107").Value = -6803.8572
# Row 132
$ws.Range("H132").Value = 62539.547
$ws.Range("I132").Value = 78306.62
$ws.Range("K132").Value = 234919.86
$ws.Range("M132").Value = -232389.86

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 97
$ws.Range("H97").Value = 946.25
$ws.Range("J97").Value = 946.25
$ws.Range("L97").Value = 2838.75
$ws.Range("N97").Value = -3830.75
# Row 114
$ws.Range("H114").Value = 13769.429
$ws.Range("I114").Value = 999.5
$ws.Range("K114").Value = 2998.5
$ws.Range("M114").Value = 255.5
# Row 122
$ws.Range("H122").Value = 1894.6842
$ws.Range("I122").Value = 1618.6
$ws.Range("J122").Value = 1993.2858
$ws.Range("K122").Value = 14567.4
$ws.Range("L122").Value = 17939.5722
$ws.Range("M122").Value = -12117.4
$ws.Range("N122").Value = -22839.5722
# Row 132
$ws.Range("H132").Value = 6672456.5
$ws.Range("I132").Value = 2067.6667
$ws.Range("J132").Value = 9531195
$ws.Range("K132").Value = 18609.0003
$ws.Range("L132").Value = 85780755
$ws.Range("M132").Value = -16079.0003
$ws.Range("N132").Value = -85785815

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 40
$ws.Range("H40").Value = 8175.75
$ws.Range("J40").Value = 8498.5
$ws.Range("L40").Value = 8498.5
$ws.Range("N40").Value = -8800.5
# Row 122
$ws.Range("H122").Value = 3888.75
$ws.Range("I122").Value = 4277.5
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 12832.5
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -10382.5
$ws.Range("N122").Value = -15400
# Row 132
$ws.Range("H132").Value = 3695.5217
$ws.Range("I132").Value = 2842.2632
$ws.Range("K132").Value = 8526.7896
$ws.Range("M132").Value = -5996.7896

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 5624.7856
$ws.Range("J40").Value = 7990
$ws.Range("L40").Value = 7990
$ws.Range("N40").Value = -8262
# Row 46
$ws.Range("H46").Value = 1203.7407
$ws.Range("J46").Value = 2601.2
$ws.Range("L46").Value = 2601.2
$ws.Range("N46").Value = -2977.2
# Row 100
$ws.Range("H100").Value = 3373.2307
$ws.Range("I100").Value = 2434.3333
$ws.Range("K100").Value = 2434.3333
$ws.Range("M100").Value = -1893.3333
# Row 109
$ws.Range("H109").Value = 56000
$ws.Range("I109").Value = 56000
$ws.Range("K109").Value = 56000
$ws.Range("M109").Value = -54613
# Row 122
$ws.Range("H122").Value = 5508.909
$ws.Range("I122").Value = 4919.8
$ws.Range("J122").Value = 5999.8335
$ws.Range("K122").Value = 14759.4
$ws.Range("L122").Value = 17999.5005
$ws.Range("M122").Value = -12309.4
$ws.Range("N122").Value = -22899.5005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 39
$ws.Range("H39").Value = 19944
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
# Row 113
$ws.Range("H113").Value = 539.7222
$ws.Range("I113").Value = 261.13333
$ws.Range("K113").Value = 783.39999
$ws.Range("M113").Value = 1386.60001
# Row 126
$ws.Range("H126").Value = 4060.5925
$ws.Range("I126").Value = 4277.68
$ws.Range("K126").Value = 12833.04
$ws.Range("M126").Value = -10363.04
# Row 132
$ws.Range("H132").Value = 3965.1719
$ws.Range("I132").Value = 4263.0376
$ws.Range("K132").Value = 12789.1128
$ws.Range("M132").Value = -10259.1128

